$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 632-633; this shifts the existing rows 632-681
# down to 634-683 and extends the sheet dimension to A1:R683.
$ws.Rows("632:633").Insert()

# New row 632: weekly data point, "Primera" quality, unidad "$/caja 36 atados"
$ws.Cells.Item(632, 1).Value2 = 6
$ws.Cells.Item(632, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(632, 3).Value2 = "Metropolitana"
$ws.Cells.Item(632, 4).Value2 = 44578
$ws.Cells.Item(632, 5).Value2 = 13
$ws.Cells.Item(632, 6).Value2 = 100112040
$ws.Cells.Item(632, 7).Value2 = "Cilantro"
$ws.Cells.Item(632, 8).Value2 = "Sin especificar"
$ws.Cells.Item(632, 9).Value2 = "Primera"
$ws.Cells.Item(632, 10).Value2 = 440
$ws.Cells.Item(632, 11).Value2 = 7500
$ws.Cells.Item(632, 12).Value2 = 8000
$ws.Cells.Item(632, 13).Value2 = 7705
$ws.Cells.Item(632, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(632, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(632, 16).Value2 = 214
$ws.Cells.Item(632, 17).Value2 = 36
$ws.Cells.Item(632, 18).Value2 = "Hortaliza"

# New row 633: weekly data point, "Primera" quality, unidad "$/docena de atados"
$ws.Cells.Item(633, 1).Value2 = 6
$ws.Cells.Item(633, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(633, 3).Value2 = "Metropolitana"
$ws.Cells.Item(633, 4).Value2 = 44578
$ws.Cells.Item(633, 5).Value2 = 13
$ws.Cells.Item(633, 6).Value2 = 100112040
$ws.Cells.Item(633, 7).Value2 = "Cilantro"
$ws.Cells.Item(633, 8).Value2 = "Sin especificar"
$ws.Cells.Item(633, 9).Value2 = "Primera"
$ws.Cells.Item(633, 10).Value2 = 380
$ws.Cells.Item(633, 11).Value2 = 14000
$ws.Cells.Item(633, 12).Value2 = 15000
$ws.Cells.Item(633, 13).Value2 = 14395
$ws.Cells.Item(633, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(633, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(633, 16).Value2 = 4798
$ws.Cells.Item(633, 17).Value2 = 3
$ws.Cells.Item(633, 18).Value2 = "Hortaliza"
